$d = $word.ActiveDocument

# --- Changes in the main document body ---
$body = $d.Content

# 1) Remove the RA number (" 000110084186 - 6 " -> two spaces)
$ok1 = $body.Find.Execute(" 000110084186 - 6 ", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "  ", 2)
Write-Host "RA number cleared:" $ok1

# 2) "QWR" (bold run right before the comma, e.g. "A QWR,") -> "TERE"
$ok2 = $body.Find.Execute("QWR", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "TERE", 2)
Write-Host "Body QWR -> TERE:" $ok2

# --- Changes in the page header ---
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdrRange = $hdr.Range

# 3) "QWER" -> "TRE"
$ok3 = $hdrRange.Find.Execute("QWER", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "TRE", 2)
Write-Host "Header QWER -> TRE:" $ok3

# 4) "QWR" -> "TERE"
$ok4 = $hdrRange.Find.Execute("QWR", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "TERE", 2)
Write-Host "Header QWR -> TERE:" $ok4

# 5) "Qwer" (5 occurrences) -> "Tre"
$ok5 = $hdrRange.Find.Execute("Qwer", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "Tre", 2)
Write-Host "Header Qwer -> Tre:" $ok5

# 6) "qwer" (3 occurrences) -> "tre"
$ok6 = $hdrRange.Find.Execute("qwer", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "tre", 2)
Write-Host "Header qwer -> tre:" $ok6
